# Add a new "entry_open" column (F) to the portfolio-mgmt sheet, with a
# default value of 0 for each existing data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "entry_open"

$ws.Range("F2").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("F4").Value = 0

$ws.Range("F5").Select() | Out-Null
